# Update the phase description in B2 to the new routine name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "FM_PHS002_Run_MLP26_to_Storeveyor"

# Add a new row (row 4) with cell A4 formatted like the existing blank
# cells in row 3 (same border/number-format/protection), so the sheet's
# used range grows to A1:B4. Copy/PasteSpecial (formats only) reproduces
# the row-3 cell formatting without disturbing its (empty) value.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
